$wb = $excel.ActiveWorkbook

# The underlying data (想去人数 / "want to go" counts) changed for four events.
# These values live identically on both the "展览" and "全部类型" worksheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 1575
    $ws.Range("F8").Value = 138
    $ws.Range("F9").Value = 62
    $ws.Range("F10").Value = 445
}
